# Weekly price-data update for the Mango / Terminal La Palmera de La Serena sheet.
# A new week of observations (Fecha = 44615, i.e. 2022-02-23) is inserted as the
# most recent entry at the top of this data block (previously starting at row 598),
# pushing the existing rows 598:627 down to 601:630 and growing the sheet from
# 627 to 630 data rows (dimension A1:T627 -> A1:T630).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above row 598 (one per quality grade: Especial/Primera/Segunda),
# shifting the existing rows 598:627 down to 601:630.
$ws.Range("A598:T600").Insert()

$qualities = @("Especial", "Primera", "Segunda")
for ($i = 0; $i -lt 3; $i++) {
    $r = 598 + $i

    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44615
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $qualities[$i]
    $ws.Cells.Item($r, 13).Value = 512
    $ws.Cells.Item($r, 14).Value = 6000
    $ws.Cells.Item($r, 15).Value = 6500
    $ws.Cells.Item($r, 16).Value = 6250
    $ws.Cells.Item($r, 17).Value = "$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = 1562
    $ws.Cells.Item($r, 20).Value = 4

    # Match the date-time number format already used by the other rows in column D.
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r + 3, 4).NumberFormat
}
